# Update profit files after running on 2025-12-15
# Appends a new data row (row 21) to the sheet with the day's results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

# Leading apostrophe forces Excel to store this as literal text instead of
# auto-parsing it into a date serial (matches the Date column's text style
# used by every other row in this sheet).
$ws.Cells.Item($row, 1).Value = "'12/15/2025"
$ws.Cells.Item($row, 2).Value = 12063.13
$ws.Cells.Item($row, 3).Value = 0.2097260978964278
$ws.Cells.Item($row, 4).Value = 0.7902739021035722
$ws.Cells.Item($row, 5).Value = -140.37
$ws.Cells.Item($row, 6).Value = -30.52
$ws.Cells.Item($row, 7).Value = -20983.58
$ws.Cells.Item($row, 8).Value = -68.76000000000001
$ws.Cells.Item($row, 9).Value = -422.9
$ws.Cells.Item($row, 10).Value = -14.32
